# CHE_convchp_fueloil.xlsx edit
#
# Semantics of the change (derived from the OOXML diff):
#   - Two new rows are inserted right after row 8 (i.e. at row 9), pushing
#     everything that used to live at row 9 onward down by two rows.
#   - The two newly inserted rows (9 and 10) are filled with what used to be
#     in row 7 and row 8 respectively (an "input_efficiency" / "constant_fxe"
#     / "fueloil" / 0.28 row, and an "output_efficiency" / "constant_fxe" /
#     "elecsupply" / 0.95 row).
#   - Row 7 and row 8 (which stay in place) get new values: the Parameter
#     (col C) becomes "input" / "output", the Type (col D) becomes
#     "configuration_fxe", and the Value (col G) becomes 1 for both.
#   - The worksheet AutoFilter range and the workbook's hidden
#     _FilterDatabase defined name both grow by two rows (L850 -> L852).
#   - The last selected cell moves to F8.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Capture the current (pre-insert) contents of rows 7 and 8 - these rows
#    will be duplicated down into the two freshly inserted rows 9 and 10.
$a7 = $ws.Range("A7").Value2
$b7 = $ws.Range("B7").Value2
$c7 = $ws.Range("C7").Value2
$d7 = $ws.Range("D7").Value2
$f7 = $ws.Range("F7").Value2
$g7 = $ws.Range("G7").Value2

$a8 = $ws.Range("A8").Value2
$b8 = $ws.Range("B8").Value2
$c8 = $ws.Range("C8").Value2
$d8 = $ws.Range("D8").Value2
$f8 = $ws.Range("F8").Value2
$g8 = $ws.Range("G8").Value2

# 2) Insert two new blank rows at row 9 - everything from the old row 9
#    onward (all the yearly data tables) shifts down by two rows.
$ws.Rows("9:10").Insert()

# 3) Populate the two new rows with the old row 7 / row 8 data.
$ws.Range("A9").Value = $a7
$ws.Range("B9").Value = $b7
$ws.Range("C9").Value = $c7
$ws.Range("D9").Value = $d7
$ws.Range("F9").Value = $f7
$ws.Range("G9").Value = $g7

$ws.Range("A10").Value = $a8
$ws.Range("B10").Value = $b8
$ws.Range("C10").Value = $c8
$ws.Range("D10").Value = $d8
$ws.Range("F10").Value = $f8
$ws.Range("G10").Value = $g8

# 4) Rewrite row 7 / row 8 in place with the new "input"/"output"
#    configuration_fxe values. Order matters for shared-string allocation
#    order (input=71, output=72, configuration_fxe=73).
$ws.Range("C7").Value = "input"
$ws.Range("C8").Value = "output"
$ws.Range("D7").Value = "configuration_fxe"
$ws.Range("D8").Value = "configuration_fxe"
$ws.Range("G7").Value = 1
$ws.Range("G8").Value = 1

# 5) Grow the AutoFilter range and the hidden _FilterDatabase name to match
#    the two extra rows now in the sheet.
$ws.AutoFilterMode = $false
$ws.Range("A5:L852").AutoFilter()
$wb.Names.Item("Sheet1!_FilterDatabase").RefersTo = "=Sheet1!`$A`$5:`$L`$852"

# 6) Leave the selection where the author last left it.
$ws.Activate()
$ws.Range("F8").Select()
